$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows (and their formatting) entirely so we start clean
$ws.Rows("1:5").Delete()

$ws.Range("A1").Value = "Mitglied"
$ws.Range("B1").Value = "Kind"
$ws.Range("C1").Value = "Klasse"
$ws.Range("D1").Value = "Eltern2"
$ws.Range("E1").Value = "Gezahlt"
$ws.Range("F1").Value = "Gezahlt Details"
$ws.Range("G1").Value = "Anzahl Kinder"
$ws.Range("H1").Value = "Zahlungszeilen"
$ws.Range("A2").Value = "Adri Mali"
$ws.Range("B2").Value = "Nik Vakl Mali"
$ws.Range("C2").Value = "Freiham Klasse 4b"
$ws.Range("D2").Value = "Adri Mali"
$ws.Range("E2").Value = 530
$ws.Range("G2").Value = 2
$ws.Range("A3").Value = "Adri Mali"
$ws.Range("B3").Value = "Alex Mali"
$ws.Range("C3").Value = "Freiham Klasse 9b"
$ws.Range("D3").Value = "Adri Mali"
$ws.Range("E3").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("A4").Value = "Max  Musterman "
$ws.Range("B4").Value = "Teo Musterman "
$ws.Range("C4").Value = "Freiham PUG 2"
$ws.Range("D4").Value = "Max  Musterman "
$ws.Range("E4").Value = 540
$ws.Range("F4").Value = "540.00 (180.00+180.00+180.00)"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = "339, 621, 647"
$ws.Range("A5").Value = "Alb Trif"
$ws.Range("B5").Value = "Valentin Trifonov"
$ws.Range("C5").Value = "Neubiberg Klasse 4a"
$ws.Range("D5").Value = "Alb Trif"
$ws.Range("E5").Value = 360
$ws.Range("F5").Value = "360.00 (180.00+180.00)"
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = "338, 517"
$ws.Range("A6").Value = "Alex Nikolov"
$ws.Range("B6").Value = "Geo Nikolov"
$ws.Range("C6").Value = "Neubiberg Klasse 5a"
$ws.Range("D6").Value = "Alex Nikolov"
$ws.Range("E6").Value = 530
$ws.Range("F6").Value = "530.00 (265.00+265.00)"
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = "337, 588"
$ws.Range("A7").Value = "Alex Nikolov"
$ws.Range("B7").Value = "Anna Nikolov"
$ws.Range("C7").Value = "Neubiberg PUG 4"
$ws.Range("D7").Value = "Alex Nikolov"
$ws.Range("E7").Value = 530
$ws.Range("F7").Value = "530.00 (265.00+265.00)"
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = "337, 588"
$ws.Range("A8").Value = "Aleks Hadzh"
$ws.Range("B8").Value = "Adri Hadzh"
$ws.Range("C8").Value = "Neubiberg Klasse 1a"
$ws.Range("D8").Value = "Aleks Hadzh"
$ws.Range("E8").Value = 530
$ws.Range("G8").Value = 1
$ws.Range("A9").Value = "Aleks Hadzh"
$ws.Range("B9").Value = "Sim Hadzh"
$ws.Range("C9").Value = "Neubiberg Klasse 1a"
$ws.Range("D9").Value = "Aleks Hadzh"
$ws.Range("E9").Value = 0
$ws.Range("G9").Value = 1

# Header row formatting
$hdr = $ws.Range("A1:H1")
$hdr.Font.Name = "Calibri"
$hdr.Font.Size = 11
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.WrapText = $false

# Column widths
$ws.Columns.Item(1).ColumnWidth = 36.666666666666664
$ws.Columns.Item(2).ColumnWidth = 32.166666666666664
$ws.Columns.Item(3).ColumnWidth = 25.833333333333336
$ws.Columns.Item(4).ColumnWidth = 27.166666666666668
$ws.Columns.Item(5).ColumnWidth = 13.0
$ws.Columns.Item(6).ColumnWidth = 32.166666666666664
$ws.Columns.Item(7).ColumnWidth = 16.666666666666668
$ws.Columns.Item(8).ColumnWidth = 16.5

# Selection
$ws.Range("F9").Select()
